$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.728.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.33%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.505.92"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.02%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.503.88"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.02%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.487"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.21%  "
$ws.Range("E10").Value = "  -0.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.57"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.430"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.08%  "
$ws.Range("E13").Value = "  -4.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.100.85"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.507.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.98%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.810.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.55%  "
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.51"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "449.17"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.626"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.650.08"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.93%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("E27").Value = "  -8.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.71"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.94"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.51"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.51%  "
$ws.Range("E31").Value = "  -3.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.169"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.45%  "
$ws.Range("E33").Value = "  +0.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.20"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.80%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.502.52"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.84"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.03"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.28"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "178.98"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.53%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0906"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.41"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.40%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "31.31"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.898"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.94"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.03%  "
$ws.Range("E51").Value = "  -1.90%  "
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.63"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.92%  "
$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.51"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -10.46%  "
